$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between I2 and J2
$ws.Range("I2").Value = 7.0
$ws.Range("J2").Value = 6.0

# Swap values between I8 and I9
$ws.Range("I8").Value = "F"
$ws.Range("I9").Value = 8.0
